# issue #5: stock data from json to db
#
# The 股票 (stock) sheet gains a new "category" column (value "normal")
# right after "property_category", plus two trailing columns,
# "source_file" (value "tmp15ba1") and "index" (value matching the row's
# existing id, 75), at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Make room for "category" in column I by shifting the existing
# date / legislator_name / legislator_id columns one place to the right
# (I:K -> J:L). Inserting (rather than overwriting) preserves each
# existing cell's formatting/style along with its value.
$ws.Columns("I:I").Insert()

# Make room at the end of the table for the two new trailing columns.
$ws.Columns("M:N").Insert()

# Header row (row 1).
$ws.Range("I1").Value2 = "category"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# Data row (row 2).
$ws.Range("I2").Value2 = "normal"
$ws.Range("M2").Value2 = "tmp15ba1"
$ws.Range("N2").Value2 = 75
